$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.620.78"
$ws.Range("E2").Value = "  -0.40%  "
$ws.Range("D3").Value = "1.891.01"
$ws.Range("E3").Value = "  -0.86%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.014"
$ws.Range("E4").Value = "  -1.33%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.42"
$ws.Range("E5").Value = "  -1.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.014"
$ws.Range("E6").Value = "  -1.60%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5138"
$ws.Range("E7").Value = "  -1.26%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3980"
$ws.Range("E8").Value = "  +0.48%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08479"
$ws.Range("E9").Value = "  +0.93%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.116"
$ws.Range("E10").Value = "  -1.90%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.292"
$ws.Range("E11").Value = "  -0.25%  "
$ws.Range("D12").Value = "1.885.10"
$ws.Range("E12").Value = "  -1.40%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.66"
$ws.Range("E13").Value = "  -0.14%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.295"
$ws.Range("E14").Value = "  -0.48%  "
$ws.Range("E15").Value = "  -0.89%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001113"
$ws.Range("E16").Value = "  -0.10%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "91.53"
$ws.Range("E17").Value = "  -0.05%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06768"
$ws.Range("E18").Value = "  -0.30%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.83"
$ws.Range("E19").Value = "  -1.13%  "
$ws.Range("E20").Value = "  -1.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.987"
$ws.Range("E21").Value = "  -2.12%  "
$ws.Range("D22").Value = "28.598.95"
$ws.Range("E22").Value = "  -0.65%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.20"
$ws.Range("E23").Value = "  -0.73%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.279"
$ws.Range("E24").Value = "  -0.68%  "
$ws.Range("D25").Value = "2.096.63"
$ws.Range("E25").Value = "  -1.57%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "161.53"
$ws.Range("E26").Value = "  -1.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.91"
$ws.Range("E27").Value = "  -0.84%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.398"
$ws.Range("E28").Value = "  -3.20%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "127.15"
$ws.Range("E29").Value = "  -0.54%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.1061"
$ws.Range("E30").Value = "  -1.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.054"
$ws.Range("E31").Value = "  -0.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.826"
$ws.Range("E32").Value = "  -2.30%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.624"
$ws.Range("E33").Value = "  -1.71%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.02443"
$ws.Range("E34").Value = "  -1.42%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.06513"
$ws.Range("E35").Value = "  -2.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.2188"
$ws.Range("E36").Value = "  -2.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.958"
$ws.Range("E37").Value = "  -5.17%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.266"
$ws.Range("E38").Value = "  +0.17%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.200"
$ws.Range("E39").Value = "  +0.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6466"
$ws.Range("E40").Value = "  -1.70%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.063"
$ws.Range("E41").Value = "  +0.66%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.28"
$ws.Range("E42").Value = "  +1.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.016"
$ws.Range("E43").Value = "  -1.37%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6091"
$ws.Range("E44").Value = "  -1.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.11"
$ws.Range("E45").Value = "  -1.21%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.723"
$ws.Range("E46").Value = "  -1.14%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.010"
$ws.Range("E47").Value = "  -0.69%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.204"
$ws.Range("E48").Value = "  -7.77%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "122.95"
$ws.Range("E49").Value = "  +0.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.213"
$ws.Range("E50").Value = "  -2.28%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06864"
$ws.Range("E51").Value = "  -1.51%  "
